$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.563.45"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "1.866.20"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4599"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07858"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9736"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.874.05"
$ws.Range("E12").Value = "  +4.52%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.681"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06938"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009996"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "28.570.64"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.109"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Value = "2.116.89"
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.778"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09315"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.268"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.320"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05775"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02064"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.704"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5613"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.751"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07168"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.138"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.138"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.830"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.413"
$ws.Range("D50").Style = "Normal"
